{"js": "// Word JS API (Office.js) script \u2014 body of `async (context) => { ... }`.\n//\n// Change described by the diff:\n//   1. The bullet \"O grupo tamb\u00e9m com a etapa anterior desenvolveu melhores\n//      capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem,\n//      tendo assim maior no\u00e7\u00e3o da estrutura das pastas.\" is reworded to\n//      \"O grupo tamb\u00e9m conseguiu com a etapa anterior desenvolver melhores\n//      capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem,\n//      tendo assim maior no\u00e7\u00e3o da sua estrutura e organiza\u00e7\u00e3o.\"\n//   2. Two new bullet items are appended right after it, in the same\n//      numbered list (numId 10 / ilvl 0, justified).\n\nconst body = context.document.body;\n\nconst oldText =\n  \"O grupo tamb\u00e9m com a etapa anterior desenvolveu melhores capacidades de \" +\n  \"programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem, tendo assim maior \" +\n  \"no\u00e7\u00e3o da estrutura das pastas.\";\n\nconst newText =\n  \"O grupo tamb\u00e9m conseguiu com a etapa anterior desenvolver melhores \" +\n  \"capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem, tendo \" +\n  \"assim maior no\u00e7\u00e3o da sua estrutura e organiza\u00e7\u00e3o.\";\n\nconst newBullet1 =\n  \"Com a ajuda de algumas dicas do professor Ricardo foi poss\u00edvel \" +\n  \"melhorar o nosso desempenho no que toca aos m\u00e9todos de organiza\u00e7\u00e3o do \" +\n  \"Trello.\";\n\nconst newBullet2 =\n  \"Com a implementa\u00e7\u00e3o do dinamismo nas nossas p\u00e1ginas est\u00e1ticas \" +\n  \"conseguimos encontrar tamb\u00e9m algumas lacunas na nossa base de dados \" +\n  \"que posteriormente foram corrigidas.\";\n\nconst results = body.search(oldText, { matchCase: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found: \" + oldText);\n}\n\nconst target = results.items[0];\n\n// Rewrite the sentence in place (keeps the paragraph's existing list /\n// alignment formatting since only the range's text is replaced).\ntarget.insertText(newText, \"Replace\");\nawait context.sync();\n\n// Re-resolve the owning paragraph so we can append the two new list items\n// right after it.\nconst paragraphs = target.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst paragraph = paragraphs.items[0];\n\nconst p1 = paragraph.insertParagraph(newBullet1, \"After\");\np1.insertParagraph(newBullet2, \"After\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n#\n# Change described by the diff:\n#   1. The bullet \"O grupo tamb\u00e9m com a etapa anterior desenvolveu melhores\n#      capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem,\n#      tendo assim maior no\u00e7\u00e3o da estrutura das pastas.\" is reworded to\n#      \"O grupo tamb\u00e9m conseguiu com a etapa anterior desenvolver melhores\n#      capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem,\n#      tendo assim maior no\u00e7\u00e3o da sua estrutura e organiza\u00e7\u00e3o.\"\n#   2. Two new bullet items are appended right after it, in the same\n#      numbered list (numId 10 / ilvl 0, justified).\n\n$d = $word.ActiveDocument\n\n$oldText = \"O grupo tamb\u00e9m com a etapa anterior desenvolveu melhores capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem, tendo assim maior no\u00e7\u00e3o da estrutura das pastas.\"\n$newText = \"O grupo tamb\u00e9m conseguiu com a etapa anterior desenvolver melhores capacidades de programa\u00e7\u00e3o e l\u00f3gica sobre esta mesma linguagem, tendo assim maior no\u00e7\u00e3o da sua estrutura e organiza\u00e7\u00e3o.\"\n$bullet1 = \"Com a ajuda de algumas dicas do professor Ricardo foi poss\u00edvel melhorar o nosso desempenho no que toca aos m\u00e9todos de organiza\u00e7\u00e3o do Trello.\"\n$bullet2 = \"Com a implementa\u00e7\u00e3o do dinamismo nas nossas p\u00e1ginas est\u00e1ticas conseguimos encontrar tamb\u00e9m algumas lacunas na nossa base de dados que posteriormente foram corrigidas.\"\n\n# Locate the target sentence and rewrite it in place (only the range's\n# text changes, so the paragraph keeps its existing numbering / alignment).\n$rng = $d.Content\n$found = $rng.Find.Execute($oldText)\nif (-not $found) {\n    throw \"Target paragraph text not found: $oldText\"\n}\n$rng.Text = $newText\n\n# The (collapsed) find range now spans the replaced sentence; its owning\n# paragraph is the bullet we just edited.\n$targetPara = $rng.Paragraphs(1)\n\n# Append the two new bullet items right after it, inheriting the list\n# numbering/justification the same way pressing Enter at end-of-paragraph\n# would in Word.\n$targetPara.Range.InsertParagraphAfter()\n$bullet1Para = $targetPara.Next()\n$bullet1Para.Range.Text = $bullet1\n\n$bullet1Para.Range.InsertParagraphAfter()\n$bullet2Para = $bullet1Para.Next()\n$bullet2Para.Range.Text = $bullet2\n"}
